$d = $word.ActiveDocument

# 1. Update both "Date de création" / "Date de version" fields: 05/07/2018 -> 06/07/2018
#    (ReplaceAll replaces every match in the Content range in one call)
$d.Content.Find.Execute("05/07/2018", $false, $false, $false, $false, $false, $true, 1, $false, "06/07/2018", 2)

# 2. Swap the order of the two "eCollection" bullet items:
#    "Traitement des retours de Mai suite à la MAJ de Mai" and
#    "Intégration des maquettes liées à la partie import"
#    Use a temporary unique marker so the two single-replace operations
#    don't clobber each other.
$d.Content.Find.Execute("Traitement des retours de Mai suite à la MAJ de Mai", $true, $false, $false, $false, $false, $true, 1, $false, "__TMP_SWAP_MARKER__", 2)
$d.Content.Find.Execute("Intégration des maquettes liées à la partie import", $true, $false, $false, $false, $false, $true, 1, $false, "Traitement des retours de Mai suite à la MAJ de Mai", 2)
$d.Content.Find.Execute("__TMP_SWAP_MARKER__", $true, $false, $false, $false, $false, $true, 1, $false, "Intégration des maquettes liées à la partie import", 2)

# 3. "Montant total : 8900€" -> "Montant total : 8900,00€"
$d.Content.Find.Execute("8900", $true, $false, $false, $false, $false, $true, 1, $false, "8900,00", 2)
